$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-11-03 Friday" "2023-11-04 Saturday"

Replace-Text "54×38=" "93×24="
Replace-Text "88×40=" "33×19="
Replace-Text "55×69=" "82×89="
Replace-Text "21×45=" "56×65="
Replace-Text "15×85=" "62×30="

Replace-Text "56×44=" "38×73="
Replace-Text "38×90=" "52×11="
Replace-Text "97×72=" "80×27="
Replace-Text "74×35=" "18×73="
Replace-Text "20×65=" "64×92="

Replace-Text "89×37=" "67×66="
Replace-Text "86×91=" "94×53="
Replace-Text "53×44=" "16×90="
Replace-Text "64×13=" "96×59="
Replace-Text "64×12=" "98×34="

Replace-Text "29×88=" "71×76="
Replace-Text "87×38=" "74×51="
Replace-Text "34×20=" "76×69="
Replace-Text "35×64=" "33×56="
Replace-Text "41×13=" "55×18="

Replace-Text "57×11=" "74×40="
Replace-Text "11×98=" "27×89="
Replace-Text "66×25=" "58×73="
Replace-Text "62×42=" "52×44="
Replace-Text "59×43=" "29×38="
